$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 15; existing rows 15-30 shift down to 16-31.
$ws.Rows("15:15").Insert()

# Populate the new row 15 with the new weekly record (same dimension
# metadata as the surrounding rows, new date/volume/price figures).
$ws.Cells.Item(15, 1).Value  = 6
$ws.Cells.Item(15, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(15, 3).Value  = "Metropolitana"
$ws.Cells.Item(15, 4).Value  = 44771
$ws.Cells.Item(15, 5).Value  = 13
$ws.Cells.Item(15, 6).Value  = 100112035
$ws.Cells.Item(15, 7).Value  = "Bruselas (repollito)"
$ws.Cells.Item(15, 8).Value  = "Sin especificar"
$ws.Cells.Item(15, 9).Value  = "Primera"
$ws.Cells.Item(15, 10).Value = 180
$ws.Cells.Item(15, 11).Value = 18000
$ws.Cells.Item(15, 12).Value = 20000
$ws.Cells.Item(15, 13).Value = 18889
$ws.Cells.Item(15, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(15, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(15, 16).Value = 1259
$ws.Cells.Item(15, 17).Value = 15
$ws.Cells.Item(15, 18).Value = "Hortaliza"
